$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 5 blank rows at row 12 to make room for the new "Docentes
#    responsaveis" block (shifts the old rows 12-22 down to 17-27).
$ws.Rows("12:16").Insert()

# 2) Fix row 10 (Objetivos:) - B/C held the wrong text before; put in the
#    actual course-objectives paragraph.
$ws.Range("B10").Value = 'Fornecer conceitos relacionados ao comportamento dos sólidos deformáveis, capacitando ao cálculo de tensões, deformações e deslocamentos em estruturas compostas por barras em regime elástico-linear sob carregamento axial, torção e flexão.Desenvolver aplicações práticas para dimensionamento de barras em condições de carregamentos mistos.Prover o conhecimento dos fenômenos de flambagem, com aplicações práticas para dimensionamento de colunas.Descrever a metodologia para análise dos estados planos de tensão e deformação, bem como a aplicação da lei de Hooke para casos multiaxiais.Apresentar conceitos básicos sobre energia de deformação.'
$ws.Range("C10").Value = 'Fornecer conceitos relacionados ao comportamento dos sólidos deformáveis, capacitando ao cálculo de tensões, deformações e deslocamentos em estruturas compostas por barras em regime elástico-linear sob carregamento axial, torção e flexão.Desenvolver aplicações práticas para dimensionamento de barras em condições de carregamentos mistos.Prover o conhecimento dos fenômenos de flambagem, com aplicações práticas para dimensionamento de colunas.Descrever a metodologia para análise dos estados planos de tensão e deformação, bem como a aplicação da lei de Hooke para casos multiaxiais.Apresentar conceitos básicos sobre energia de deformação.'

# 3) New row 12: "Docentes responsaveis:" label (column A only).
$ws.Range("A12").Value = 'Docentes responsáveis:'

# 4) New rows 13-16: the four professor names in B/C (no column A).
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '471420 - Carlos Antonio Reis Pereira Baptista'
$ws.Range("C13").Value = '471420 - Carlos Antonio Reis Pereira Baptista'
$ws.Range("A14").Clear()
$ws.Range("B14").Value = '3480026 - João Paulo Pascon'
$ws.Range("C14").Value = '3480026 - João Paulo Pascon'
$ws.Range("A15").Clear()
$ws.Range("B15").Value = '5840793 - Sérgio Schneider'
$ws.Range("C15").Value = '5840793 - Sérgio Schneider'
$ws.Range("A16").Clear()
$ws.Range("B16").Value = '7797767 - Viktor Pastoukhov'
$ws.Range("C16").Value = '7797767 - Viktor Pastoukhov'

# 5) Match column B/C formatting (wrap text / font) on the newly-created
#    rows to the rest of the table by copying it from existing cells.
$ws.Range("B10").Copy()
$ws.Range("B13:B16").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13:C16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 6) Row 17 (was 12): "Programa resumido:" now holds the real short-syllabus text.
$ws.Range("B17").Value = 'Considerações fundamentais. Conceito de tensão. Conceito de deformação. Lei de Hooke. Carga Axial. Torção em barras de seção circular. Flexão em vigas isostáticas de seção simétrica. Cargas combinadas. Flambagem de colunas. Análise de Tensão e Deformação. Lei de Hooke Multiaxial. Energia de deformação.'
$ws.Range("C17").Value = 'Considerações fundamentais. Conceito de tensão. Conceito de deformação. Lei de Hooke. Carga Axial. Torção em barras de seção circular. Flexão em vigas isostáticas de seção simétrica. Cargas combinadas. Flambagem de colunas. Análise de Tensão e Deformação. Lei de Hooke Multiaxial. Energia de deformação.'

# 7) Row 19 (was 14): "Programa:" now holds the full syllabus text.
$ws.Range("B19").Value = '1. Considerações fundamentais: Propósito da Mecânica dos Sólidos; Carregamentos e Esforços Solicitantes.2. Conceito de tensão: Tensão Normal; Tensão Cisalhante; Tensões admissíveis.3. Conceito de deformação: Deformação Normal; Deformação por Cisalhamento.4. Lei de Hooke: Elasticidade linear e o Módulo de Young; Lei de Hooke para Cisalhamento.5. Carga Axial: Deslocamentos em sistemas isostáticos; Efeitos da Temperatura; Sistemas Hiperestáticos.6. Torção em barras de seção circular: Momento de inércia polar; Análise das tensões em eixos de seção maciça e seção vazada; Cálculo das rotações relativas entre seções adjacentes; Eixos estaticamente indeterminados; Torção e tração combinadas.7. Flexão em vigas isostáticas de seção simétrica: Forças concentradas e forças distribuídas; Diagramas de força cortante e momento fletor para uma viga carregada; Momento de inércia, eixos principais de inércia; Flexão em Vigas de Seção Simétrica; Determinação das Tensões Normais; Deflexões em vigas: equação diferencial da linha elástica; Tensões de cisalhamento em vigas. Tensões de cisalhamento em barras de paredes finas.8. Cargas combinadas: Modos Mistos de Carregamento. Projeto de barras submetidas a cargas axiais, transversais e torcionais.9. Flambagem de colunas: Raio de giração. Fórmula de Euler para colunas biarticuladas. Fatores de correção para outras condições de contorno. Projeto de colunas de aço e de outras ligas submetidas a um carregamento centrado.10. Análise de Tensão e Deformação: Variação da Tensão com o Plano de Corte; Estado Plano de Tensão; Tensões Principais e Máxima Tensão de Cisalhamento; O Círculo de Mohr para Tensão Plana; Tensão Triaxial; Transformação do Estado Plano de Deformação.11. Lei de Hooke Multiaxial: Elasticidade, Homogeneidade e Isotropia; Coeficiente de Poisson; Lei de Hooke para Tensão Triaxial em Materiais Isotrópicos; Relações entre as Constantes Elásticas; Aplicação em Vasos de Pressão de Paredes Finas.12. Energia de deformação: Densidade de energia de deformação. Energia de deformação elástica para tensões normais. Energia de deformação elástica para tensões de cisalhamento. Projeto para carregamento por impacto. Métodos de energia: teorema de Castigliano e suas aplicações.'
$ws.Range("C19").Value = '1. Considerações fundamentais: Propósito da Mecânica dos Sólidos; Carregamentos e Esforços Solicitantes.2. Conceito de tensão: Tensão Normal; Tensão Cisalhante; Tensões admissíveis.3. Conceito de deformação: Deformação Normal; Deformação por Cisalhamento.4. Lei de Hooke: Elasticidade linear e o Módulo de Young; Lei de Hooke para Cisalhamento.5. Carga Axial: Deslocamentos em sistemas isostáticos; Efeitos da Temperatura; Sistemas Hiperestáticos.6. Torção em barras de seção circular: Momento de inércia polar; Análise das tensões em eixos de seção maciça e seção vazada; Cálculo das rotações relativas entre seções adjacentes; Eixos estaticamente indeterminados; Torção e tração combinadas.7. Flexão em vigas isostáticas de seção simétrica: Forças concentradas e forças distribuídas; Diagramas de força cortante e momento fletor para uma viga carregada; Momento de inércia, eixos principais de inércia; Flexão em Vigas de Seção Simétrica; Determinação das Tensões Normais; Deflexões em vigas: equação diferencial da linha elástica; Tensões de cisalhamento em vigas. Tensões de cisalhamento em barras de paredes finas.8. Cargas combinadas: Modos Mistos de Carregamento. Projeto de barras submetidas a cargas axiais, transversais e torcionais.9. Flambagem de colunas: Raio de giração. Fórmula de Euler para colunas biarticuladas. Fatores de correção para outras condições de contorno. Projeto de colunas de aço e de outras ligas submetidas a um carregamento centrado.10. Análise de Tensão e Deformação: Variação da Tensão com o Plano de Corte; Estado Plano de Tensão; Tensões Principais e Máxima Tensão de Cisalhamento; O Círculo de Mohr para Tensão Plana; Tensão Triaxial; Transformação do Estado Plano de Deformação.11. Lei de Hooke Multiaxial: Elasticidade, Homogeneidade e Isotropia; Coeficiente de Poisson; Lei de Hooke para Tensão Triaxial em Materiais Isotrópicos; Relações entre as Constantes Elásticas; Aplicação em Vasos de Pressão de Paredes Finas.12. Energia de deformação: Densidade de energia de deformação. Energia de deformação elástica para tensões normais. Energia de deformação elástica para tensões de cisalhamento. Projeto para carregamento por impacto. Métodos de energia: teorema de Castigliano e suas aplicações.'

# 8) Row 22 (was 17): "Metodo:" now holds the evaluation-method text.
$ws.Range("B22").Value = 'Os alunos serão avaliados por meio de três conjuntos de notas: duas provas escritas (P1 e P2) envolvendo o conteúdo teórico ministrado em sala de aula; exercícios (EX) propostos regularmente para serem entregues e discutidos na aula subsequente; e seminários (SM) em grupo ao final da disciplina.'
$ws.Range("C22").Value = 'Os alunos serão avaliados por meio de três conjuntos de notas: duas provas escritas (P1 e P2) envolvendo o conteúdo teórico ministrado em sala de aula; exercícios (EX) propostos regularmente para serem entregues e discutidos na aula subsequente; e seminários (SM) em grupo ao final da disciplina.'

# 9) Row 23 (was 18): "Criterio:" now holds the final-grade formula text.
$ws.Range("B23").Value = 'Nota Final (NF) = 70%((P1+P2)/2)+20%(EX)+10%(SM).'
$ws.Range("C23").Value = 'Nota Final (NF) = 70%((P1+P2)/2)+20%(EX)+10%(SM).'

# 10) Row 24 (was 19): "Norma de recuperacao:" now holds the make-up exam rule.
$ws.Range("B24").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'
$ws.Range("C24").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'

# 11) Row 25 (was 20): "Bibliografia:" now holds the full reading list.
$ws.Range("B25").Value = '1. J.M. GERE. Mecânica dos Materiais. São Paulo: Pioneira Thomson Learning, 2003, 698p.2. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF. Resistência dos Materiais. São Paulo: McGraw Hill. 4a Ed., 2006, 758p.3. R.R. CRAIG, Jr. Mecânica dos Materiais. Rio de Janeiro LTC. 2a Ed., 2003, 552p.4. R.C. HIBBELER. Resistência dos Materiais. São Paulo: Pearson Prentice Hall. 5a Ed., 2006, 670p.5. A.C. UGURAL. Mecânica dos Materiais. Rio de Janeiro LTC, 2009, 638p.6. A.R. RAGAB, S.E. BAYOUMI. Engineering Solid Mechanics, Fundamentals and Applications. New York: CRC Press, 1999, 921p.7. POPOV, E. P. Introdução à Mecânica dos Sólidos, São Paulo: Edgard Blücher, 1978, 552p.8. A. HIGDON, E.H. OHLSEN, W.B. STILES, J.A. WEESE, W.F. RILEY. Mecânica dos Materiais. Rio de Janeiro: Guanabara Dois. 3a Ed., 1981, 549p.'
$ws.Range("C25").Value = '1. J.M. GERE. Mecânica dos Materiais. São Paulo: Pioneira Thomson Learning, 2003, 698p.2. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF. Resistência dos Materiais. São Paulo: McGraw Hill. 4a Ed., 2006, 758p.3. R.R. CRAIG, Jr. Mecânica dos Materiais. Rio de Janeiro LTC. 2a Ed., 2003, 552p.4. R.C. HIBBELER. Resistência dos Materiais. São Paulo: Pearson Prentice Hall. 5a Ed., 2006, 670p.5. A.C. UGURAL. Mecânica dos Materiais. Rio de Janeiro LTC, 2009, 638p.6. A.R. RAGAB, S.E. BAYOUMI. Engineering Solid Mechanics, Fundamentals and Applications. New York: CRC Press, 1999, 921p.7. POPOV, E. P. Introdução à Mecânica dos Sólidos, São Paulo: Edgard Blücher, 1978, 552p.8. A. HIGDON, E.H. OHLSEN, W.B. STILES, J.A. WEESE, W.F. RILEY. Mecânica dos Materiais. Rio de Janeiro: Guanabara Dois. 3a Ed., 1981, 549p.'

Write-Output "done"
